$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Badge Magnet" line item (row 7). This shifts the Total row
# up to row 7, tightens the shared-formula range, the SUM formula, and
# drops the now-unused shared strings ("Badge Magnet" and its URL).
$ws.Rows("7:7").Delete()

# Update the sheet's cell selection to match the author's final cursor
# position.
$ws.Range("G13").Select()
